$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": update version/date, replace the broken "Contact" rows
#     with a single "Jurisdiction" row, and fill in the publisher value.
$meta = $wb.Worksheets.Item("Metadata")

# Version bump
$meta.Cells.Item(3, 2).Value = "6.0.0"

# Date bump
$meta.Cells.Item(8, 2).Value = "2022-01-21T20:46:54+00:00"

# Publisher value (was blank)
$meta.Cells.Item(9, 2).Value = "Alvearie Team"

# Turn the old "Contact" / "No display for ContactDetail" row into
# "Jurisdiction" / "United States of America"
$meta.Cells.Item(10, 1).Value = "Jurisdiction"
$meta.Cells.Item(10, 2).Value = "United States of America"

# Remove the duplicate "Contact" row that followed it (row 11); everything
# below shifts up by one, turning the old A1:B21 range into A1:B20.
$meta.Rows.Item(11).Delete()

# --- Sheet "Elements": the "Extension" root row's Short/Definition text
#     was regenerated.
$elements = $wb.Worksheets.Item("Elements")
$elements.Cells.Item(2, 11).Value = "Evaluated Ouptut"
$elements.Cells.Item(2, 12).Value = "Attachment for content created as output when producing the insight."
